# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for the rows whose values changed in this data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.556.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.078.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.626'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.47'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +5.64%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.387'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0765'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.386.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.16'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.13%  '
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.066.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '37.694.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +17.11%  '
$ws.Range('E21').Value = '  +3.08%  '
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '226.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.17%  '
$ws.Range('E25').Value = '  +2.99%  '
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.62%  '
$ws.Range('E28').Value = '  +9.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.06'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.03%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0627'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('E35').Value = '  +6.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.90%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.56%  '
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +19.56%  '
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.474.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').Value = '  +6.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '95.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.87%  '
$ws.Range('E47').Value = '  +4.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('E49').Value = '  +4.02%  '
$ws.Range('E50').Value = '  +5.54%  '
$ws.Range('E51').Value = '  +1.73%  '
